$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "0.964 (0.964 ± 0.000)"
$ws.Range("C2").Value = "00:04:02 (00:04:33 ± 00:00:12)"
$ws.Range("D2").Value = "00:00:08 (00:00:13 ± 00:00:03)"
$ws.Range("B3").Value = "0.966 (0.956 ± 0.005)"
$ws.Range("C3").Value = "00:00:38 (00:01:20 ± 00:00:21)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B4").Value = "0.944 (0.918 ± 0.010)"
$ws.Range("C4").Value = "00:01:40 (00:02:04 ± 00:00:23)"
$ws.Range("D4").Value = "00:00:01 (00:00:01 ± 00:00:00)"
$ws.Range("B5").Value = "0.964 (0.953 ± 0.007)"
$ws.Range("C5").Value = "00:05:07 (00:05:14 ± 00:00:04)"
$ws.Range("D5").Value = "00:00:01 (00:00:02 ± 00:00:01)"
$ws.Range("B6").Value = "0.975 (0.958 ± 0.006)"
$ws.Range("C6").Value = "00:04:57 (00:05:01 ± 00:00:02)"
$ws.Range("D6").Value = "00:00:02 (00:00:06 ± 00:00:03)"
$ws.Range("B7").Value = "0.959 (0.948 ± 0.007)"
$ws.Range("C7").Value = "00:05:00 (00:05:04 ± 00:00:03)"
$ws.Range("D7").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B8").Value = "0.968 (0.956 ± 0.007)"
$ws.Range("C8").Value = "00:04:51 (00:06:19 ± 00:01:49)"
$ws.Range("D8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B9").Value = "0.971 (0.958 ± 0.005)"
$ws.Range("C9").Value = "00:05:00 (00:05:02 ± 00:00:02)"
$ws.Range("D9").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B10").Value = "0.970 (0.954 ± 0.007)"
$ws.Range("C10").Value = "00:04:29 (00:04:29 ± 00:00:00)"
$ws.Range("D10").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B11").Value = "0.939 (0.888 ± 0.028)"
$ws.Range("C11").Value = "00:05:06 (00:05:06 ± 00:00:00)"
$ws.Range("D11").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B12").Value = "0.244 (0.227 ± 0.013)"
$ws.Range("C12").Value = "00:02:13 (00:02:38 ± 00:00:10)"
$ws.Range("D12").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B13").Value = "0.966 (0.948 ± 0.008)"
$ws.Range("C13").Value = "00:00:07 (00:00:10 ± 00:00:02)"
$ws.Range("D13").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B14").Value = "0.968 (0.953 ± 0.006)"
$ws.Range("C14").Value = "00:00:27 (00:00:29 ± 00:00:00)"
$ws.Range("D14").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B15").Value = "0.964 (0.958 ± 0.003)"
$ws.Range("C15").Value = "00:02:46 (00:04:33 ± 00:00:40)"
$ws.Range("D15").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B16").Value = "0.969 (0.956 ± 0.006)"
$ws.Range("C16").Value = "00:00:34 (00:00:37 ± 00:00:02)"
$ws.Range("D16").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B17").Value = "0.971 (0.954 ± 0.006)"
$ws.Range("C17").Value = "00:05:02 (00:05:26 ± 00:00:16)"
$ws.Range("D17").Value = "00:00:00 (00:00:00 ± 00:00:00)"
